# Updates cryptocurrency price (D) and volume change (E) columns
# to match the latest scraped data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.849.92'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.39%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.645.58'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.20%  '
$ws.Range('E4').Value = '  +1.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.10'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.73%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.502'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.18%  '
$ws.Range('E7').Value = '  +1.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.252'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.23%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0625'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.49%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.07'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.95%  '
$ws.Range('E11').Value = '  -0.15%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.869.32'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.46%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.648.86'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.75%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.19'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.06%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.527'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.90%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.96'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.53%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.838.39'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.22%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0738'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.67%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '213.21'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.84%  '
$ws.Range('E20').Value = '  +0.98%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.34'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.39%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.41'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +14.85%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.28'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.60%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.39'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.91%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.59'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.64%  '
$ws.Range('E26').Value = '  +1.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.118'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.35%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.17'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.33%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.67'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.77%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0520'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.26%  '
$ws.Range('E31').Value = '  +0.36%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.32'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.70%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.98'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.39%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.278.03'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.53'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.97%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.45'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.02%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0177'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.09%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.537'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.20%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.827'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.54%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.02'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.12%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.815'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.14%  '
$ws.Range('E42').Value = '  -1.73%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.38'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.35%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.795.05'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.43%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '91.75'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.86%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '58.75'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.49%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.61'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.37%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0₆0103'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.67%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0520'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.72%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.65'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.11%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0976'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.15%  '
